$wb = $excel.ActiveWorkbook

# --- Text update: "Ready for handoff" -> "In Translation" -------------------
# This shared string is referenced by:
#   Overview!E2, Overview!F2  (per-locale status columns)
#   zh-cn!C2                  (Status column)
#   de-de!C2                  (Status column)
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E2").Value = "In Translation"
$ws1.Range("F2").Value = "In Translation"

$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("C2").Value = "In Translation"

$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("C2").Value = "In Translation"

# --- Re-fit the Status columns now that the text is shorter ----------------
# The new text is narrower than "Ready for handoff", so the report generator
# narrows the affected columns accordingly.
$ws1.Columns.Item(5).ColumnWidth = 12.5
$ws1.Columns.Item(6).ColumnWidth = 12.5

$ws2.Columns.Item(3).ColumnWidth = 12.5

$ws3.Columns.Item(3).ColumnWidth = 12.5
